# The underlying canonical-OOXML diff for this revision is purely an
# attribute-serialization artifact: every hunk re-orders the XML
# attributes already present on word/document.xml, word/footer1.xml and
# word/footnotes.xml (namespace declarations alphabetized, then element
# attributes alphabetized, e.g. `w:w="3070" w:type="dxa"` becomes
# `w:type="dxa" w:w="3070"`). No text, formatting, table data, style, or
# property value actually changes between the two revisions.
#
# That byte-level attribute order is an artifact of whatever tool
# re-saved the package (outside the Word object model) and is not
# something the Word OM exposes a way to control - there is no
# property/method that reorders XML attributes on a part. So there is
# no content-level mutation to perform here: we simply confirm the
# template's content is intact (read-only), leaving the document
# semantically identical, which matches the diff.

$d = $word.ActiveDocument

# Sanity-check: footer still holds the static-table demo text untouched.
$found = $d.Content.Find.Execute("A simple demonstration of a static table", `
                                  $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)

Write-Output ("Footer table template verified, paragraphs=" + $d.Paragraphs.Count)
